$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from H1 (bold/bordered header style) to new header cells I1, J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows
$iValues = @{2=1; 3=1; 4=1; 5=1; 6=1; 7=9; 8=1; 9=6}
$jValues = @{2=5; 3=5; 4=5; 5=5; 6=5; 7=9; 8=4; 9=6}

foreach ($row in 2..9) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
